# corrected data cleaning for pre/post/total fixation data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Strip the header row's special formatting (bold font + border + center/top alignment). ---
$ws.Range("A1:P1").ClearFormats()

# --- 2. The old "Unnamed: 0" header label is removed (cell becomes blank). ---
$ws.Range("A1").Value = ""

# --- 3. Insert a new row above the old "Revisit count" row and populate it with the
#        newly-computed "Respondent ratio (%)" data. This pushes the old rows 3-8 down
#        to rows 4-9 and extends the used range to row 11 (new blank trailing row). ---
$ws.Rows.Item(3).Insert()

$ws.Range("A3").Value = "Respondent ratio (%)"
$ws.Range("B3").Value = 100
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 100
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 100
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 100
$ws.Range("J3").Value = 100
$ws.Range("K3").Value = 100
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = 0
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = 0

# --- 4. Re-computed values for the metrics rows that shifted down (rows 4-9). ---
# Revisit count (row 4)
$ws.Range("B4").Value = 8
$ws.Range("G4").Value = 1
$ws.Range("J4").Value = 2
$ws.Range("K4").Value = 24

# Fixation count (row 5)
$ws.Range("B5").Value = 13
$ws.Range("G5").Value = 2
$ws.Range("I5").Value = 2
$ws.Range("J5").Value = 3
$ws.Range("K5").Value = 99

# Dwell time (ms) (row 6)
$ws.Range("B6").Value = 2636.23
$ws.Range("G6").Value = 433.77
$ws.Range("J6").Value = 567.1799999999999
$ws.Range("K6").Value = 32091.68

# Dwell time (%) (row 7)
$ws.Range("B7").Value = 2.53
$ws.Range("G7").Value = 0.42
$ws.Range("I7").Value = 0.61
$ws.Range("J7").Value = 0.54
$ws.Range("K7").Value = 30.77

# Fixation duration (ms) (row 8)
$ws.Range("B8").Value = 202.79
$ws.Range("G8").Value = 216.89
$ws.Range("J8").Value = 189.06
$ws.Range("K8").Value = 324.16

# First fixation duration (ms) (row 9) - values unchanged from before the shift,
# nothing further required here.
